$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 12097.2
$ws.Range("J62").Value = 5361.2
$ws.Range("L62").Value = 5361.2
$ws.Range("N62").Value = -6609.2
$ws.Range("H65").Value = 12097.2
$ws.Range("J65").Value = 5361.2
$ws.Range("L65").Value = 26806
$ws.Range("N65").Value = -33046
$ws.Range("H137").Value = 24433.105
$ws.Range("I137").Value = 33243.16
$ws.Range("J137").Value = 7363.625
$ws.Range("K137").Value = 99729.48000000001
$ws.Range("L137").Value = 22090.875
$ws.Range("M137").Value = -97179.48000000001
$ws.Range("N137").Value = -27190.875
$ws.Range("H138").Value = 1702.0941
$ws.Range("I138").Value = 847.4186
$ws.Range("J138").Value = 2577.1191
$ws.Range("K138").Value = 2542.2558
$ws.Range("L138").Value = 7731.3573
$ws.Range("M138").Value = 2597.7442
$ws.Range("N138").Value = -18011.3573
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1129542.8
$ws.Range("I32").Value = 1321788.1
$ws.Range("J32").Value = 5646.5386
$ws.Range("K32").Value = 1321788.1
$ws.Range("L32").Value = 5646.5386
$ws.Range("M32").Value = -1321501.1
$ws.Range("N32").Value = -6220.5386
$ws.Range("H60").Value = 250027260
$ws.Range("I60").Value = 1000000000
$ws.Range("K60").Value = 1000000000
$ws.Range("M60").Value = -999999267
$ws.Range("H61").Value = 1262.5769
$ws.Range("I61").Value = 836.9231
$ws.Range("J61").Value = 2539.5386
$ws.Range("K61").Value = 836.9231
$ws.Range("L61").Value = 2539.5386
$ws.Range("M61").Value = -624.9231
$ws.Range("N61").Value = -2963.5386
$ws.Range("H74").Value = 25211.62
$ws.Range("I74").Value = 40006.92
$ws.Range("J74").Value = 1169.25
$ws.Range("K74").Value = 40006.92
$ws.Range("L74").Value = 1169.25
$ws.Range("M74").Value = -39132.92
$ws.Range("N74").Value = -2917.25
$ws.Range("H77").Value = 25211.62
$ws.Range("I77").Value = 40006.92
$ws.Range("J77").Value = 1169.25
$ws.Range("K77").Value = 200034.6
$ws.Range("L77").Value = 5846.25
$ws.Range("M77").Value = -195666.6
$ws.Range("N77").Value = -14582.25
$ws.Range("H102").Value = 1316.5
$ws.Range("I102").Value = 1279.8
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1279.8
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 342.2
$ws.Range("N102").Value = -4744
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 1069.3334
$ws.Range("I122").Value = 804
$ws.Range("K122").Value = 2412
$ws.Range("M122").Value = 38
$ws.Range("H123").Value = 34014.5
$ws.Range("J123").Value = 34014.5
$ws.Range("L123").Value = 34014.5
$ws.Range("N123").Value = -43814.5
$ws.Range("H132").Value = 1840472.8
$ws.Range("I132").Value = 2220627
$ws.Range("J132").Value = 674666.75
$ws.Range("K132").Value = 6661881
$ws.Range("L132").Value = 2024000.25
$ws.Range("M132").Value = -6659351
$ws.Range("N132").Value = -2029060.25
$ws.Range("H136").Value = 1262.5769
$ws.Range("I136").Value = 836.9231
$ws.Range("J136").Value = 2539.5386
$ws.Range("K136").Value = 2510.7693
$ws.Range("L136").Value = 7618.6158
$ws.Range("M136").Value = 39.23070000000007
$ws.Range("N136").Value = -12718.6158
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 43976.75
$ws.Range("J61").Value = 43976.75
$ws.Range("L61").Value = 43976.75
$ws.Range("N61").Value = -44602.75
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 34740
$ws.Range("J123").Value = 34740
$ws.Range("L123").Value = 34740
$ws.Range("N123").Value = -44540
$ws.Range("H134").Value = 17884.016
$ws.Range("I134").Value = 1040.5283
$ws.Range("J134").Value = 86553.62
$ws.Range("K134").Value = 3121.5849
$ws.Range("L134").Value = 259660.86
$ws.Range("M134").Value = -586.5848999999998
$ws.Range("N134").Value = -264730.86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9254.078
$ws.Range("I31").Value = 6991.6587
$ws.Range("J31").Value = 18530
$ws.Range("K31").Value = 6991.6587
$ws.Range("L31").Value = 18530
$ws.Range("M31").Value = -6696.6587
$ws.Range("N31").Value = -19120
$ws.Range("H34").Value = 9254.078
$ws.Range("I34").Value = 6991.6587
$ws.Range("J34").Value = 18530
$ws.Range("K34").Value = 6991.6587
$ws.Range("L34").Value = 18530
$ws.Range("M34").Value = -6789.6587
$ws.Range("N34").Value = -18934
$ws.Range("H51").Value = 9612.625
$ws.Range("J51").Value = 9612.625
$ws.Range("L51").Value = 9612.625
$ws.Range("N51").Value = -11084.625
$ws.Range("H59").Value = 12509.889
$ws.Range("J59").Value = 12509.889
$ws.Range("L59").Value = 12509.889
$ws.Range("N59").Value = -14799.889
$ws.Range("H60").Value = 9221.444
$ws.Range("J60").Value = 9221.444
$ws.Range("L60").Value = 9221.444
$ws.Range("N60").Value = -10243.444
$ws.Range("H61").Value = 9612.625
$ws.Range("J61").Value = 9612.625
$ws.Range("L61").Value = 9612.625
$ws.Range("N61").Value = -10308.625
$ws.Range("H124").Value = 19297.5
$ws.Range("J124").Value = 19297.5
$ws.Range("L124").Value = 19297.5
$ws.Range("N124").Value = -24207.5
$ws.Range("H125").Value = 24103.5
$ws.Range("J125").Value = 24103.5
$ws.Range("L125").Value = 24103.5
$ws.Range("N125").Value = -29023.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5310.1904
$ws.Range("I5").Value = 532.4545000000001
$ws.Range("J5").Value = 10565.7
$ws.Range("K5").Value = 1597.3635
$ws.Range("L5").Value = 31697.1
$ws.Range("M5").Value = -1485.3635
$ws.Range("N5").Value = -31921.1
$ws.Range("H122").Value = 502.1282
$ws.Range("J122").Value = 545.9655
$ws.Range("L122").Value = 4913.6895
$ws.Range("N122").Value = -9813.6895
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 13500
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -23380
$ws.Range("H131").Value = 30488650
$ws.Range("I131").Value = 510
$ws.Range("J131").Value = 34723110
$ws.Range("K131").Value = 1530
$ws.Range("L131").Value = 104169330
$ws.Range("M131").Value = 3510
$ws.Range("N131").Value = -104179410
$ws.Range("H135").Value = 5310.1904
$ws.Range("I135").Value = 532.4545000000001
$ws.Range("J135").Value = 10565.7
$ws.Range("K135").Value = 4792.0905
$ws.Range("L135").Value = 95091.3
$ws.Range("M135").Value = -2257.0905
$ws.Range("N135").Value = -100161.3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 35000
$ws.Range("J124").Value = 35000
$ws.Range("L124").Value = 35000
$ws.Range("N124").Value = -44820
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -44920
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 25200
$ws.Range("J124").Value = 25200
$ws.Range("L124").Value = 25200
$ws.Range("N124").Value = -35020
$ws.Range("H125").Value = 50715
$ws.Range("J125").Value = 50715
$ws.Range("L125").Value = 50715
$ws.Range("N125").Value = -60555
$ws.Range("H132").Value = 234073.78
$ws.Range("I132").Value = 67008.39
$ws.Range("J132").Value = 480694.16
$ws.Range("K132").Value = 201025.17
$ws.Range("L132").Value = 1442082.48
$ws.Range("M132").Value = -198495.17
$ws.Range("N132").Value = -1447142.48
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 3117.7874
$ws.Range("I132").Value = 707.7857
$ws.Range("J132").Value = 6669.3687
$ws.Range("K132").Value = 2123.3571
$ws.Range("L132").Value = 20008.1061
$ws.Range("M132").Value = 406.6428999999998
$ws.Range("N132").Value = -25068.1061
$ws.Range("H136").Value = 947711.6
$ws.Range("I136").Value = 1083119
$ws.Range("J136").Value = 500867.44
$ws.Range("K136").Value = 3249357
$ws.Range("L136").Value = 1502602.32
$ws.Range("M136").Value = -3246807
$ws.Range("N136").Value = -1507702.32
